$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and 1h volume change (column E) values
# Column D values must stay as text (they use "." as thousands separators, e.g. 27.282.32),
# so force text number format before assignment, then restore the default "Normal" style
# so the cell keeps looking like the untouched cells around it.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.282.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.900.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5216"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.76%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3779"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.66%  "

$ws.Range("E9").Value = "  +1.18%  "

$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9026"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08161"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.900.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.29%  "

$ws.Range("E15").Value = "  +1.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008633"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.309.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.099"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.87%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.427"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.304"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.55%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.746"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.833"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.928"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09251"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05068"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7990"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.98%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.448"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.951"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.605"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5710"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02005"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.81%  "

$ws.Range("E40").Value = "  +0.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.011"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.581"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1520"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4896"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.45%  "

$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.623"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.85%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05950"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.48%  "
